# Insert a new weekly price observation as row 328 in the "Cebolla" data
# sheet. This shifts the existing rows 328-368 down to 329-369 (the sheet
# dimension grows from A1:R368 to A1:R369) and fills the new row 328 with
# the new observation's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 328..368 down by one row, leaving row 328 empty for the new
# data point.
$ws.Rows.Item(328).Insert()

$ws.Range("A328").Value = 7
$ws.Range("B328").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C328").Value = "Ñuble"
$ws.Range("D328").Value = 44505
$ws.Range("E328").Value = 16
$ws.Range("F328").Value = 100112004
$ws.Range("G328").Value = "Cebolla"
$ws.Range("H328").Value = "Sin especificar"
$ws.Range("I328").Value = "1a nueva(o)"
$ws.Range("J328").Value = 20000
$ws.Range("K328").Value = 950
$ws.Range("L328").Value = 1000
$ws.Range("M328").Value = 975
$ws.Range("N328").Value = "`$/paquete 10 unidades (volumen en unidades)"
$ws.Range("O328").Value = "Región de O'Higgins"
$ws.Range("P328").Value = 98
$ws.Range("Q328").Value = 10
$ws.Range("R328").Value = "Hortaliza"
